# Update TestData.xlsx row 5 (fresh QA data) for integration test setup.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TGLTestData")

# New test email address (replaces previous QA fresh email in F5)
$ws.Range("F5").Value = "amadasun@usc.edu"

# New PID value used across PIDIntegation/PIDConfirmed/PIDAssignment/PIDDeferAccept
$ws.Range("H5").Value = 4383359
$ws.Range("I5").Value = 4383359
$ws.Range("J5").Value = 4383359
$ws.Range("K5").Value = 4383359

# Restore the active selection to D5, matching the saved workbook view state
$ws.Range("D5").Select()
